$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 2423.1892  # ALC!H112 2533.75 -> 2423.1892
$ws.Cells.Item(112, 10).Value = 2696.5  # ALC!J112 2799.4285 -> 2696.5
$ws.Cells.Item(112, 12).Value = 8089.5  # ALC!L112 8398.2855 -> 8089.5
$ws.Cells.Item(112, 14).Value = -10305.5  # ALC!N112 -10614.2855 -> -10305.5

$ws.Cells.Item(132, 8).Value = 1686.14  # ALC!H132 1701.6346 -> 1686.14
$ws.Cells.Item(132, 9).Value = 1451.1538  # ALC!I132 1479.4103 -> 1451.1538
$ws.Cells.Item(132, 10).Value = 2519.2727  # ALC!J132 2368.3076 -> 2519.2727
$ws.Cells.Item(132, 11).Value = 4353.4614  # ALC!K132 4438.2309 -> 4353.4614
$ws.Cells.Item(132, 12).Value = 7557.8181  # ALC!L132 7104.9228 -> 7557.8181
$ws.Cells.Item(132, 13).Value = -1823.4614  # ALC!M132 -1908.2309 -> -1823.4614
$ws.Cells.Item(132, 14).Value = -12617.8181  # ALC!N132 -12164.9228 -> -12617.8181

$ws.Cells.Item(135, 8).Value = 604.55554  # ALC!H135 593.413 -> 604.55554
$ws.Cells.Item(135, 9).Value = 595.2069  # ALC!I135 591.4828 -> 595.2069
$ws.Cells.Item(135, 10).Value = 621.5  # ALC!J135 596.7059 -> 621.5
$ws.Cells.Item(135, 11).Value = 5356.8621  # ALC!K135 5323.3452 -> 5356.8621
$ws.Cells.Item(135, 12).Value = 5593.5  # ALC!L135 5370.3531 -> 5593.5
$ws.Cells.Item(135, 13).Value = -2821.8621  # ALC!M135 -2788.3452 -> -2821.8621
$ws.Cells.Item(135, 14).Value = -10663.5  # ALC!N135 -10440.3531 -> -10663.5

$ws.Cells.Item(137, 8).Value = 2716.4119  # ALC!H137 2801.8064 -> 2716.4119
$ws.Cells.Item(137, 9).Value = 1763.7273  # ALC!I137 1869.8 -> 1763.7273
$ws.Cells.Item(137, 10).Value = 3172.0435  # ALC!J137 3245.6191 -> 3172.0435
$ws.Cells.Item(137, 11).Value = 5291.1819  # ALC!K137 5609.4 -> 5291.1819
$ws.Cells.Item(137, 12).Value = 9516.130500000001  # ALC!L137 9736.8573 -> 9516.130500000001
$ws.Cells.Item(137, 13).Value = -2741.1819  # ALC!M137 -3059.4 -> -2741.1819
$ws.Cells.Item(137, 14).Value = -14616.1305  # ALC!N137 -14836.8573 -> -14616.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7471.3096  # ARM!H32 7483.976 -> 7471.3096
$ws.Cells.Item(32, 9).Value = 6709.606  # ARM!I32 6500.853 -> 6709.606
$ws.Cells.Item(32, 10).Value = 10264.223  # ARM!J32 11662.25 -> 10264.223
$ws.Cells.Item(32, 11).Value = 6709.606  # ARM!K32 6500.853 -> 6709.606
$ws.Cells.Item(32, 12).Value = 10264.223  # ARM!L32 11662.25 -> 10264.223
$ws.Cells.Item(32, 13).Value = -6422.606  # ARM!M32 -6213.853 -> -6422.606
$ws.Cells.Item(32, 14).Value = -10838.223  # ARM!N32 -12236.25 -> -10838.223

$ws.Cells.Item(88, 8).Value = 2401.5715  # ARM!H88 2434.3333 -> 2401.5715
$ws.Cells.Item(88, 9).Value = 2005.5  # ARM!I88 2006 -> 2005.5
$ws.Cells.Item(88, 10).Value = 2560  # ARM!J88 2520 -> 2560
$ws.Cells.Item(88, 11).Value = 2005.5  # ARM!K88 2006 -> 2005.5
$ws.Cells.Item(88, 12).Value = 2560  # ARM!L88 2520 -> 2560
$ws.Cells.Item(88, 13).Value = -1599.5  # ARM!M88 -1600 -> -1599.5
$ws.Cells.Item(88, 14).Value = -3372  # ARM!N88 -3332 -> -3372

$ws.Cells.Item(91, 8).Value = 2401.5715  # ARM!H91 2434.3333 -> 2401.5715
$ws.Cells.Item(91, 9).Value = 2005.5  # ARM!I91 2006 -> 2005.5
$ws.Cells.Item(91, 10).Value = 2560  # ARM!J91 2520 -> 2560
$ws.Cells.Item(91, 11).Value = 2005.5  # ARM!K91 2006 -> 2005.5
$ws.Cells.Item(91, 12).Value = 2560  # ARM!L91 2520 -> 2560
$ws.Cells.Item(91, 13).Value = -601.5  # ARM!M91 -602 -> -601.5
$ws.Cells.Item(91, 14).Value = -5368  # ARM!N91 -5328 -> -5368

$ws.Cells.Item(123, 8).Value = 0  # ARM!H123 25341.1 -> 0
$ws.Cells.Item(123, 10).Value = 0  # ARM!J123 25341.1 -> 0
$ws.Cells.Item(123, 12).Value = 0  # ARM!L123 25341.1 -> 0
$ws.Cells.Item(123, 14).ClearContents()  # ARM!N123 -35141.1 -> (removed)

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 251372.5  # BSM!H20 84316.664 -> 251372.5
$ws.Cells.Item(20, 9).Value = 334533.34  # BSM!I20 200799.4 -> 334533.34
$ws.Cells.Item(20, 10).Value = 1890  # BSM!J20 1114.7142 -> 1890
$ws.Cells.Item(20, 11).Value = 334533.34  # BSM!K20 200799.4 -> 334533.34
$ws.Cells.Item(20, 12).Value = 1890  # BSM!L20 1114.7142 -> 1890
$ws.Cells.Item(20, 13).Value = -334286.34  # BSM!M20 -200552.4 -> -334286.34
$ws.Cells.Item(20, 14).Value = -2384  # BSM!N20 -1608.7142 -> -2384

$ws.Cells.Item(105, 8).Value = 6805804  # BSM!H105 5104923 -> 6805804
$ws.Cells.Item(105, 9).Value = 8406239  # BSM!I105 6496183 -> 8406239
$ws.Cells.Item(105, 10).Value = 3955.25  # BSM!J105 3636.8333 -> 3955.25
$ws.Cells.Item(105, 11).Value = 8406239  # BSM!K105 6496183 -> 8406239
$ws.Cells.Item(105, 12).Value = 3955.25  # BSM!L105 3636.8333 -> 3955.25
$ws.Cells.Item(105, 13).Value = -8404492  # BSM!M105 -6494436 -> -8404492
$ws.Cells.Item(105, 14).Value = -7449.25  # BSM!N105 -7130.8333 -> -7449.25

$ws.Cells.Item(107, 8).Value = 26108.957  # BSM!H107 22894.111 -> 26108.957
$ws.Cells.Item(107, 9).Value = 29125.35  # BSM!I107 26851.318 -> 29125.35
$ws.Cells.Item(107, 10).Value = 5999.6665  # BSM!J107 5482.4 -> 5999.6665
$ws.Cells.Item(107, 11).Value = 29125.35  # BSM!K107 26851.318 -> 29125.35
$ws.Cells.Item(107, 12).Value = 5999.6665  # BSM!L107 5482.4 -> 5999.6665
$ws.Cells.Item(107, 13).Value = -27205.35  # BSM!M107 -24931.318 -> -27205.35
$ws.Cells.Item(107, 14).Value = -9839.6665  # BSM!N107 -9322.4 -> -9839.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(35, 8).Value = 439  # CRP!H35 325 -> 439
$ws.Cells.Item(35, 9).Value = 439  # CRP!I35 325 -> 439
$ws.Cells.Item(35, 11).Value = 439  # CRP!K35 325 -> 439
$ws.Cells.Item(35, 13).Value = -145  # CRP!M35 -31 -> -145

$ws.Cells.Item(62, 8).Value = 58311.668  # CRP!H62 74172.14 -> 58311.668
$ws.Cells.Item(62, 9).Value = 65200.625  # CRP!I62 86000.836 -> 65200.625
$ws.Cells.Item(62, 11).Value = 65200.625  # CRP!K62 86000.836 -> 65200.625
$ws.Cells.Item(62, 13).Value = -64576.625  # CRP!M62 -85376.836 -> -64576.625

$ws.Cells.Item(65, 8).Value = 58311.668  # CRP!H65 74172.14 -> 58311.668
$ws.Cells.Item(65, 9).Value = 65200.625  # CRP!I65 86000.836 -> 65200.625
$ws.Cells.Item(65, 11).Value = 326003.125  # CRP!K65 430004.18 -> 326003.125
$ws.Cells.Item(65, 13).Value = -322883.125  # CRP!M65 -426884.18 -> -322883.125

$ws.Cells.Item(122, 8).Value = 2055.6  # CRP!H122 1719.5 -> 2055.6
$ws.Cells.Item(122, 9).Value = 1750  # CRP!I122 1500 -> 1750
$ws.Cells.Item(122, 10).Value = 2259.3333  # CRP!J122 1792.6666 -> 2259.3333
$ws.Cells.Item(122, 11).Value = 5250  # CRP!K122 4500 -> 5250
$ws.Cells.Item(122, 12).Value = 6777.999899999999  # CRP!L122 5377.9998 -> 6777.999899999999
$ws.Cells.Item(122, 13).Value = -2800  # CRP!M122 -2050 -> -2800
$ws.Cells.Item(122, 14).Value = -11677.9999  # CRP!N122 -10277.9998 -> -11677.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(20, 8).Value = 800  # CUL!H20 1000 -> 800
$ws.Cells.Item(20, 9).Value = 800  # CUL!I20 400 -> 800
$ws.Cells.Item(20, 10).Value = 0  # CUL!J20 4000 -> 0
$ws.Cells.Item(20, 11).Value = 2400  # CUL!K20 1200 -> 2400
$ws.Cells.Item(20, 12).Value = 0  # CUL!L20 12000 -> 0
$ws.Cells.Item(20, 13).Value = -2173  # CUL!M20 -973 -> -2173
$ws.Cells.Item(20, 14).ClearContents()  # CUL!N20 -12454 -> (removed)

$ws.Cells.Item(21, 8).Value = 750  # CUL!H21 111112710 -> 750
$ws.Cells.Item(21, 9).Value = 750  # CUL!I21 800 -> 750
$ws.Cells.Item(21, 10).Value = 0  # CUL!J21 166668670 -> 0
$ws.Cells.Item(21, 11).Value = 2250  # CUL!K21 2400 -> 2250
$ws.Cells.Item(21, 12).Value = 0  # CUL!L21 500006010 -> 0
$ws.Cells.Item(21, 13).Value = -2077  # CUL!M21 -2227 -> -2077
$ws.Cells.Item(21, 14).ClearContents()  # CUL!N21 -500006356 -> (removed)

$ws.Cells.Item(39, 8).Value = 4354.0713  # CUL!H39 4524.467 -> 4354.0713
$ws.Cells.Item(39, 9).Value = 650  # CUL!I39 0 -> 650
$ws.Cells.Item(39, 10).Value = 4971.4165  # CUL!J39 4524.467 -> 4971.4165
$ws.Cells.Item(39, 11).Value = 1950  # CUL!K39 0 -> 1950
$ws.Cells.Item(39, 12).Value = 14914.2495  # CUL!L39 13573.401 -> 14914.2495
$ws.Cells.Item(39, 13).Value = -1656  # CUL!M39 (blank) -> -1656
$ws.Cells.Item(39, 14).Value = -15502.2495  # CUL!N39 -14161.401 -> -15502.2495

$ws.Cells.Item(51, 8).Value = 1821.75  # CUL!H51 1850.75 -> 1821.75
$ws.Cells.Item(51, 9).Value = 514.8  # CUL!I51 604 -> 514.8
$ws.Cells.Item(51, 10).Value = 4000  # CUL!J51 2266.3333 -> 4000
$ws.Cells.Item(51, 11).Value = 1544.4  # CUL!K51 1812 -> 1544.4
$ws.Cells.Item(51, 12).Value = 12000  # CUL!L51 6798.999899999999 -> 12000
$ws.Cells.Item(51, 13).Value = -1084.4  # CUL!M51 -1352 -> -1084.4
$ws.Cells.Item(51, 14).Value = -12920  # CUL!N51 -7718.999899999999 -> -12920

$ws.Cells.Item(59, 8).Value = 0  # CUL!H59 1600 -> 0
$ws.Cells.Item(59, 9).Value = 0  # CUL!I59 1600 -> 0
$ws.Cells.Item(59, 11).Value = 0  # CUL!K59 4800 -> 0
$ws.Cells.Item(59, 13).ClearContents()  # CUL!M59 -4260 -> (removed)

$ws.Cells.Item(63, 8).Value = 265751.5  # CUL!H63 303473.16 -> 265751.5
$ws.Cells.Item(63, 9).Value = 422402.4  # CUL!I63 702170.7 -> 422402.4
$ws.Cells.Item(63, 10).Value = 4666.6665  # CUL!J63 4450 -> 4666.6665
$ws.Cells.Item(63, 11).Value = 1267207.2  # CUL!K63 2106512.1 -> 1267207.2
$ws.Cells.Item(63, 12).Value = 13999.9995  # CUL!L63 13350 -> 13999.9995
$ws.Cells.Item(63, 13).Value = -1266458.2  # CUL!M63 -2105763.1 -> -1266458.2
$ws.Cells.Item(63, 14).Value = -15497.9995  # CUL!N63 -14848 -> -15497.9995

$ws.Cells.Item(66, 8).Value = 265751.5  # CUL!H66 303473.16 -> 265751.5
$ws.Cells.Item(66, 9).Value = 422402.4  # CUL!I66 702170.7 -> 422402.4
$ws.Cells.Item(66, 10).Value = 4666.6665  # CUL!J66 4450 -> 4666.6665
$ws.Cells.Item(66, 11).Value = 3801621.6  # CUL!K66 6319536.3 -> 3801621.6
$ws.Cells.Item(66, 12).Value = 41999.9985  # CUL!L66 40050 -> 41999.9985
$ws.Cells.Item(66, 13).Value = -3797877.6  # CUL!M66 -6315792.3 -> -3797877.6
$ws.Cells.Item(66, 14).Value = -49487.9985  # CUL!N66 -47538 -> -49487.9985

$ws.Cells.Item(68, 8).Value = 1328.4  # CUL!H68 1274.093 -> 1328.4
$ws.Cells.Item(68, 9).Value = 1058.0238  # CUL!I68 991.21277 -> 1058.0238
$ws.Cells.Item(68, 10).Value = 1627.2368  # CUL!J68 1615 -> 1627.2368
$ws.Cells.Item(68, 11).Value = 3174.0714  # CUL!K68 2973.63831 -> 3174.0714
$ws.Cells.Item(68, 12).Value = 4881.7104  # CUL!L68 4845 -> 4881.7104
$ws.Cells.Item(68, 13).Value = -2363.0714  # CUL!M68 -2162.63831 -> -2363.0714
$ws.Cells.Item(68, 14).Value = -6503.7104  # CUL!N68 -6467 -> -6503.7104

$ws.Cells.Item(71, 8).Value = 1328.4  # CUL!H71 1274.093 -> 1328.4
$ws.Cells.Item(71, 9).Value = 1058.0238  # CUL!I71 991.21277 -> 1058.0238
$ws.Cells.Item(71, 10).Value = 1627.2368  # CUL!J71 1615 -> 1627.2368
$ws.Cells.Item(71, 11).Value = 9522.214199999999  # CUL!K71 8920.914929999999 -> 9522.214199999999
$ws.Cells.Item(71, 12).Value = 14645.1312  # CUL!L71 14535 -> 14645.1312
$ws.Cells.Item(71, 13).Value = -5466.214199999999  # CUL!M71 -4864.914929999999 -> -5466.214199999999
$ws.Cells.Item(71, 14).Value = -22757.1312  # CUL!N71 -22647 -> -22757.1312

$ws.Cells.Item(107, 8).Value = 1320.1428  # CUL!H107 1286.027 -> 1320.1428
$ws.Cells.Item(107, 9).Value = 1135.2727  # CUL!I107 1077.25 -> 1135.2727
$ws.Cells.Item(107, 10).Value = 1633  # CUL!J107 1671.4615 -> 1633
$ws.Cells.Item(107, 11).Value = 3405.8181  # CUL!K107 3231.75 -> 3405.8181
$ws.Cells.Item(107, 12).Value = 4899  # CUL!L107 5014.3845 -> 4899
$ws.Cells.Item(107, 13).Value = -1485.8181  # CUL!M107 -1311.75 -> -1485.8181
$ws.Cells.Item(107, 14).Value = -8739  # CUL!N107 -8854.3845 -> -8739

$ws.Cells.Item(122, 8).Value = 593.64703  # CUL!H122 628.63635 -> 593.64703
$ws.Cells.Item(122, 9).Value = 477.25  # CUL!I122 450.05554 -> 477.25
$ws.Cells.Item(122, 10).Value = 759.9286  # CUL!J122 842.93335 -> 759.9286
$ws.Cells.Item(122, 11).Value = 4295.25  # CUL!K122 4050.49986 -> 4295.25
$ws.Cells.Item(122, 12).Value = 6839.3574  # CUL!L122 7586.40015 -> 6839.3574
$ws.Cells.Item(122, 13).Value = -1845.25  # CUL!M122 -1600.49986 -> -1845.25
$ws.Cells.Item(122, 14).Value = -11739.3574  # CUL!N122 -12486.40015 -> -11739.3574

$ws.Cells.Item(131, 8).Value = 20836640  # CUL!H131 16669496 -> 20836640
$ws.Cells.Item(131, 9).Value = 1461.25  # CUL!I131 1293.2142 -> 1461.25
$ws.Cells.Item(131, 10).Value = 27781700  # CUL!J131 21742428 -> 27781700
$ws.Cells.Item(131, 11).Value = 4383.75  # CUL!K131 3879.6426 -> 4383.75
$ws.Cells.Item(131, 12).Value = 83345100  # CUL!L131 65227284 -> 83345100
$ws.Cells.Item(131, 13).Value = 656.25  # CUL!M131 1160.3574 -> 656.25
$ws.Cells.Item(131, 14).Value = -83355180  # CUL!N131 -65237364 -> -83355180

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 261392  # GSM!H15 175863.83 -> 261392
$ws.Cells.Item(15, 10).Value = 261392  # GSM!J15 175863.83 -> 261392
$ws.Cells.Item(15, 12).Value = 261392  # GSM!L15 175863.83 -> 261392
$ws.Cells.Item(15, 14).Value = -261968  # GSM!N15 -176439.83 -> -261968

$ws.Cells.Item(81, 8).Value = 261392  # GSM!H81 175863.83 -> 261392
$ws.Cells.Item(81, 10).Value = 261392  # GSM!J81 175863.83 -> 261392
$ws.Cells.Item(81, 12).Value = 261392  # GSM!L81 175863.83 -> 261392
$ws.Cells.Item(81, 14).Value = -263388  # GSM!N81 -177859.83 -> -263388

$ws.Cells.Item(84, 8).Value = 261392  # GSM!H84 175863.83 -> 261392
$ws.Cells.Item(84, 10).Value = 261392  # GSM!J84 175863.83 -> 261392
$ws.Cells.Item(84, 12).Value = 784176  # GSM!L84 527591.49 -> 784176
$ws.Cells.Item(84, 14).Value = -794160  # GSM!N84 -537575.49 -> -794160

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 20033334  # LTW!H20 15037500 -> 20033334
$ws.Cells.Item(20, 9).Value = 20033334  # LTW!I20 15037500 -> 20033334
$ws.Cells.Item(20, 11).Value = 20033334  # LTW!K20 15037500 -> 20033334
$ws.Cells.Item(20, 13).Value = -20033108  # LTW!M20 -15037274 -> -20033108

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 1525000  # WVR!H20 5025000 -> 1525000
$ws.Cells.Item(20, 9).Value = 1525000  # WVR!I20 5025000 -> 1525000
$ws.Cells.Item(20, 11).Value = 1525000  # WVR!K20 5025000 -> 1525000
$ws.Cells.Item(20, 13).Value = -1524760  # WVR!M20 -5024760 -> -1524760

$ws.Cells.Item(109, 8).Value = 24932.777  # WVR!H109 27832.889 -> 24932.777
$ws.Cells.Item(109, 10).Value = 24932.777  # WVR!J109 27832.889 -> 24932.777
$ws.Cells.Item(109, 12).Value = 24932.777  # WVR!L109 27832.889 -> 24932.777
$ws.Cells.Item(109, 14).Value = -27706.777  # WVR!N109 -30606.889 -> -27706.777
